$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIO - Flour mill")

# Delete the entire row for "Manildra Flour Mill" (row 3), shifting rows below it up.
$ws.Rows.Item(3).Delete()

# Restore the selection to match the post-edit state.
$ws.Range("H21").Select()
